# Auto-generated edit script: update market-data columns (H-N) per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 5: Met a Sticky End / Animal Glue
$ws.Range("H5").Value = 638
$ws.Range("I5").Value = 63.333332
$ws.Range("J5").Value = 1500
$ws.Range("K5").Value = 63.333332
$ws.Range("L5").Value = 1500
$ws.Range("M5").Value = 51.666668
$ws.Range("N5").Value = -1730

# Row 29: Dripping with Venom / Weak Blinding Potion
$ws.Range("H29").Value = 916.6667
$ws.Range("J29").Value = 1325
$ws.Range("L29").Value = 3975
$ws.Range("N29").Value = -4537

# Row 33: Glazed and Confused / Clear Glass Lens
$ws.Range("H33").Value = 389.83334
$ws.Range("I33").Value = 335.90475
$ws.Range("K33").Value = 335.90475
$ws.Range("M33").Value = -106.90475

# Row 38: Just Give Him a Serum / Hi-Potion of Strength
$ws.Range("H38").Value = 168.66667
$ws.Range("I38").Value = 168.66667
$ws.Range("K38").Value = 506.00001
$ws.Range("M38").Value = -134.00001

# Row 70: Consecrating Congregation / Holy Water
$ws.Range("H70").Value = 1343.909
$ws.Range("I70").Value = 1140.5714
$ws.Range("J70").Value = 1699.75
$ws.Range("K70").Value = 3421.7142
$ws.Range("L70").Value = 5099.25
$ws.Range("M70").Value = -3151.7142
$ws.Range("N70").Value = -5639.25

# Row 73: Curbing the Contagion (L) / Holy Water
$ws.Range("H73").Value = 1343.909
$ws.Range("I73").Value = 1140.5714
$ws.Range("J73").Value = 1699.75
$ws.Range("K73").Value = 3421.7142
$ws.Range("L73").Value = 5099.25
$ws.Range("M73").Value = -2485.7142
$ws.Range("N73").Value = -6971.25

# Row 74: Adhesive of Antipathy / Wing Glue
$ws.Range("H74").Value = 4466878.5
$ws.Range("I74").Value = 2189.0557
$ws.Range("J74").Value = 12503320
$ws.Range("K74").Value = 2189.0557
$ws.Range("L74").Value = 12503320
$ws.Range("M74").Value = -1253.0557
$ws.Range("N74").Value = -12505192

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 3706809.5
$ws.Range("I76").Value = 3233.3333
$ws.Range("J76").Value = 9262174
$ws.Range("K76").Value = 3233.3333
$ws.Range("L76").Value = 9262174
$ws.Range("M76").Value = -2918.3333
$ws.Range("N76").Value = -9262804

# Row 77: It's Gonna Grow Back (L) / Wing Glue
$ws.Range("H77").Value = 4466878.5
$ws.Range("I77").Value = 2189.0557
$ws.Range("J77").Value = 12503320
$ws.Range("K77").Value = 10945.2785
$ws.Range("L77").Value = 62516600
$ws.Range("M77").Value = -6265.2785
$ws.Range("N77").Value = -62525960

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 3706809.5
$ws.Range("I79").Value = 3233.3333
$ws.Range("J79").Value = 9262174
$ws.Range("K79").Value = 3233.3333
$ws.Range("L79").Value = 9262174
$ws.Range("M79").Value = -2141.3333
$ws.Range("N79").Value = -9264358

# Row 100: Asking for a Friend / Beetle Glue
$ws.Range("H100").Value = 2266.4285
$ws.Range("I100").Value = 1526.8182
$ws.Range("J100").Value = 3080
$ws.Range("K100").Value = 1526.8182
$ws.Range("L100").Value = 3080
$ws.Range("M100").Value = -985.8181999999999
$ws.Range("N100").Value = -4162

# Row 129: Practical Command / Commanding Craftsman's Draught
$ws.Range("H129").Value = 944.7
$ws.Range("I129").Value = 457.2
$ws.Range("J129").Value = 998.86664
$ws.Range("K129").Value = 1371.6
$ws.Range("L129").Value = 2996.59992
$ws.Range("M129").Value = 3628.4
$ws.Range("N129").Value = -12996.59992

# Row 141: Remedy for Reason / Grade 1 Gemdraught of Mind
$ws.Range("H141").Value = 2320.5557
$ws.Range("I141").Value = 1735.625
$ws.Range("J141").Value = 7000
$ws.Range("K141").Value = 5206.875
$ws.Range("L141").Value = 21000
$ws.Range("M141").Value = -26.875
$ws.Range("N141").Value = -31360

$ws = $wb.Worksheets.Item("BSM")
# Row 22: Riveting Run / Iron Rivets
$ws.Range("H22").Value = 344
$ws.Range("I22").Value = 344
$ws.Range("K22").Value = 344
$ws.Range("M22").Value = -171

# Row 134: Ruthenium Supremium / Ruthenium Ingot
$ws.Range("H134").Value = 2798.2646
$ws.Range("I134").Value = 2849.4062
$ws.Range("K134").Value = 8548.2186
$ws.Range("M134").Value = -6013.2186

# Row 138: Bladewinner / Titanium Gold Greatsword
$ws.Range("H138").Value = 39870
$ws.Range("J138").Value = 39870
$ws.Range("L138").Value = 39870
$ws.Range("N138").Value = -50150

$ws = $wb.Worksheets.Item("CRP")
# Row 58: You Do the Heavy Lifting / Mahogany Lumber
$ws.Range("H58").Value = 24994.238
$ws.Range("I58").Value = 1265.5
$ws.Range("J58").Value = 167366.67
$ws.Range("K58").Value = 1265.5
$ws.Range("L58").Value = 167366.67
$ws.Range("M58").Value = -1062.5
$ws.Range("N58").Value = -167772.67

# Row 74: License to Heal / Dark Chestnut Rod
$ws.Range("H74").Value = 30246.143
$ws.Range("J74").Value = 30246.143
$ws.Range("L74").Value = 30246.143
$ws.Range("N74").Value = -31994.143

# Row 77: Purified Polyrhythm (L) / Dark Chestnut Rod
$ws.Range("H77").Value = 30246.143
$ws.Range("J77").Value = 30246.143
$ws.Range("L77").Value = 90738.429
$ws.Range("N77").Value = -99474.429

# Row 88: Hold on Adamantite / Adamantite Spear
$ws.Range("H88").Value = 31375
$ws.Range("J88").Value = 31375
$ws.Range("L88").Value = 31375
$ws.Range("N88").Value = -32187

# Row 91: Spears for Stone Vigilantes (L) / Adamantite Spear
$ws.Range("H91").Value = 31375
$ws.Range("J91").Value = 31375
$ws.Range("L91").Value = 31375
$ws.Range("N91").Value = -34183

# Row 136: Turali Quality / Dark Mahogany Lumber
$ws.Range("H136").Value = 24994.238
$ws.Range("I136").Value = 1265.5
$ws.Range("J136").Value = 167366.67
$ws.Range("K136").Value = 3796.5
$ws.Range("L136").Value = 502100.01
$ws.Range("M136").Value = -1246.5
$ws.Range("N136").Value = -507200.01

$ws = $wb.Worksheets.Item("CUL")
# Row 11: Putting the Squeeze On / Orange Juice
$ws.Range("H11").Value = 128.83333
$ws.Range("I11").Value = 105.25
$ws.Range("J11").Value = 176
$ws.Range("K11").Value = 315.75
$ws.Range("L11").Value = 528
$ws.Range("M11").Value = -175.75
$ws.Range("N11").Value = -808

# Row 17: Chew the Fat / Grilled Dodo
$ws.Range("H17").Value = 800
$ws.Range("J17").Value = 800
$ws.Range("L17").Value = 2400
$ws.Range("N17").Value = -2738

# Row 26: A Grape Idea / Grape Juice
$ws.Range("H26").Value = 355.81818
$ws.Range("I26").Value = 120.666664
$ws.Range("J26").Value = 638
$ws.Range("K26").Value = 361.999992
$ws.Range("L26").Value = 1914
$ws.Range("M26").Value = -73.99999200000002
$ws.Range("N26").Value = -2490

# Row 34: Fever Pitch / Chamomile Tea
$ws.Range("H34").Value = 1000
$ws.Range("J34").Value = 1000
$ws.Range("L34").Value = 3000
$ws.Range("N34").Value = -3168

# Row 68: Such a Butter Face / Fermented Butter
$ws.Range("H68").Value = 1140.4546
$ws.Range("J68").Value = 1154.5
$ws.Range("L68").Value = 3463.5
$ws.Range("N68").Value = -5085.5

# Row 71: No Margarine of Error (L) / Fermented Butter
$ws.Range("H71").Value = 1140.4546
$ws.Range("J71").Value = 1154.5
$ws.Range("L71").Value = 10390.5
$ws.Range("N71").Value = -18502.5

# Row 81: It Goes Down Smoothly / Frozen Spirits
$ws.Range("H81").Value = 4363.1113
$ws.Range("J81").Value = 5436.4287
$ws.Range("L81").Value = 16309.2861
$ws.Range("N81").Value = -18555.2861

# Row 84: Quenching the Flame (L) / Frozen Spirits
$ws.Range("H84").Value = 4363.1113
$ws.Range("J84").Value = 5436.4287
$ws.Range("L84").Value = 48927.85830000001
$ws.Range("N84").Value = -60159.85830000001

# Row 122: Salt of the North / Northern Sea Salt
$ws.Range("H122").Value = 757.1667
$ws.Range("I122").Value = 384.75
$ws.Range("J122").Value = 1502
$ws.Range("K122").Value = 3462.75
$ws.Range("L122").Value = 13518
$ws.Range("M122").Value = -1012.75
$ws.Range("N122").Value = -18418

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 818.02
$ws.Range("J131").Value = 830.5361
$ws.Range("L131").Value = 2491.6083
$ws.Range("N131").Value = -12571.6083

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3409.762
$ws.Range("I80").Value = 2792.0833
$ws.Range("J80").Value = 4233.3335
$ws.Range("K80").Value = 2792.0833
$ws.Range("L80").Value = 4233.3335
$ws.Range("M80").Value = -1794.0833
$ws.Range("N80").Value = -6229.3335

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3409.762
$ws.Range("I83").Value = 2792.0833
$ws.Range("J83").Value = 4233.3335
$ws.Range("K83").Value = 13960.4165
$ws.Range("L83").Value = 21166.6675
$ws.Range("M83").Value = -8968.416499999999
$ws.Range("N83").Value = -31150.6675

# Row 132: On Board for Lar / Lar Ingot
$ws.Range("H132").Value = 20839.035
$ws.Range("I132").Value = 2928.3333
$ws.Range("J132").Value = 74571.14
$ws.Range("K132").Value = 8784.999899999999
$ws.Range("L132").Value = 223713.42
$ws.Range("M132").Value = -6254.999899999999
$ws.Range("N132").Value = -228773.42

$ws = $wb.Worksheets.Item("LTW")
# Row 94: Fitting In / Gaganaskin Hat of Aiming
$ws.Range("H94").Value = 36443.332
$ws.Range("J94").Value = 36443.332
$ws.Range("L94").Value = 36443.332
$ws.Range("N94").Value = -37795.332

# Row 104: Brace Yourselves / Gazelleskin Bracers of Fending
$ws.Range("H104").Value = 24428.625
$ws.Range("J104").Value = 24428.625
$ws.Range("L104").Value = 24428.625
$ws.Range("N104").Value = -31416.625

# Row 135: Dreams of Ja / Crocodileskin Leg Wraps of Scouting
$ws.Range("H135").Value = 25000
$ws.Range("J135").Value = 25000
$ws.Range("L135").Value = 25000
$ws.Range("N135").Value = -35140
